# Update metric data: Mon Apr 28 13:13:21 UTC 2025
# Appends the next timestamp/metric sample as a new row at the bottom
# of the data table (row 15), growing the used range from A1:B14 to
# A1:B15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 15

$ws.Cells.Item($newRow, 1).Value = "2025-04-28 13:13:21"
$ws.Cells.Item($newRow, 2).Value = 239
